# Mapeo de categorias y campos de producto, contruccion del excel
#
# Expands "Tabla1" from 8 columns (Categoria 1-5, Nombre, Precio, Imagen) to
# 17 columns, inserting product-detail fields (Marca, Dimensiones, Largo,
# Ancho, Alto, Diametro, Color) and renaming the old "Nombre"/"Precio" slots
# to "Uso"/"Contenido", keeping "Imagen" as the last column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Grow the table (and its AutoFilter) out to column Q so there is room for
# all of the new product-detail columns.
$tbl.Resize($ws.Range("A1:Q1048576"))

# Final header order (this also drives xl/sharedStrings.xml + the table's
# <tableColumns> list, since each ListColumn's Name mirrors its header
# cell text).
$headers = @(
    "Categoria 1",
    "Categoria 2",
    "Categoria 3",
    "Categoria 4",
    "Categoria 5",
    "Marca",
    "Nombre",
    "Precio",
    "Dimensiones",
    "Largo",
    "Ancho",
    "Alto",
    "Diametro",
    "Color",
    "Uso",
    "Contenido",
    "Imagen"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Match the column sizing applied alongside the new fields.
$ws.Range("A1:Q1").ColumnWidth = 10.33

# Restore the header-row selection/active cell to the new last column.
$ws.Range("Q2").Select() | Out-Null
